$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New observation records (rows 10-13) appended to the "Artfynd" sheet,
# mirroring the existing column layout (row 1 headers).

$rows = @(
    @{
        A = 112163316; B = 90666; C = "Ovaliderad"; D = "LC"; E = 4364
        F = "Dropptaggsvamp"; G = "Hydnellum ferrugineum"; H = "(Fr.:Fr.) P. Karst."
        I = "21"; J = "fruktkroppar"; K = ""; N = ""
        P = "Källdalen SV, Vstm"; Q = 558054.0685264322; R = 6628664.783651764; S = 10
        T = "Västmanland"; U = "Surahammar"; V = "Västmanland"; W = "Ramnäs"
        Y = "2023-09-17"; Z = "00:00"; AA = "2023-09-17"; AB = "00:00"
        AC = "Tätt bestånd."
        AD = $false; AE = $false; AF = ""; AG = $false
        AI = "Gammal barrblandskog"
        AT = ""
        AW = "Tom Sävström"; AX = "Tom Sävström"; AY = ""
    },
    @{
        A = 112163367; B = 90666; C = "Ovaliderad"; D = "LC"; E = 4364
        F = "Dropptaggsvamp"; G = "Hydnellum ferrugineum"; H = "(Fr.:Fr.) P. Karst."
        I = ""; J = ""; K = ""; N = ""
        P = "Källdalen SV, Vstm"; Q = 558082.638313611; R = 6628611.362187758; S = 10
        T = "Västmanland"; U = "Surahammar"; V = "Västmanland"; W = "Ramnäs"
        Y = "2023-09-17"; Z = "00:00"; AA = "2023-09-17"; AB = "00:00"
        AD = $false; AE = $false; AF = ""; AG = $false
        AI = "Gammal barrblandskog"
        AT = ""
        AW = "Tom Sävström"; AX = "Tom Sävström"; AY = ""
    },
    @{
        A = 112164162; B = 90666; C = "Ovaliderad"; D = "LC"; E = 4364
        F = "Dropptaggsvamp"; G = "Hydnellum ferrugineum"; H = "(Fr.:Fr.) P. Karst."
        I = ""; J = ""; K = ""; N = ""
        P = "Månses hål, Vstm"; Q = 558022.2854273538; R = 6628309.856821301; S = 10
        T = "Västmanland"; U = "Surahammar"; V = "Västmanland"; W = "Ramnäs"
        Y = "2023-09-17"; Z = "00:00"; AA = "2023-09-17"; AB = "00:00"
        AD = $false; AE = $false; AF = ""; AG = $false
        AI = "Gammal barrblandskog, renlavsmarker"
        AT = ""
        AW = "Tom Sävström"; AX = "Tom Sävström"; AY = ""
    },
    @{
        A = 112163866; B = 90666; C = "Ovaliderad"; D = "LC"; E = 4364
        F = "Dropptaggsvamp"; G = "Hydnellum ferrugineum"; H = "(Fr.:Fr.) P. Karst."
        I = ""; J = ""; K = ""; N = ""
        P = "Månses hål, Vstm"; Q = 558019.5269801348; R = 6628292.695551688; S = 10
        T = "Västmanland"; U = "Surahammar"; V = "Västmanland"; W = "Ramnäs"
        Y = "2023-09-17"; Z = "00:00"; AA = "2023-09-17"; AB = "00:00"
        AD = $false; AE = $false; AF = ""; AG = $false
        AI = "Gammal barrblandskog, renlavsmarker"
        AT = ""
        AW = "Tom Sävström"; AX = "Tom Sävström"; AY = ""
    }
)

# Columns whose source values are text even though they look like a
# number/date (e.g. "21", "2023-09-17") - force Text format first so
# Excel doesn't auto-convert them to a number/date serial.
$textColumns = @("I", "Y", "AA")

$startRow = 10
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    foreach ($col in $data.Keys) {
        $value = $data[$col]
        $cell = $ws.Range("$col$r")
        if ($textColumns -contains $col -and $value -ne "") {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $value
    }
}
